$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Player roster data (Oyuncu Adi, Pozisyon, Takim) for rows 2-19
$data = @(
    ,@('Cole Anthony', 'PG', 'Orlando Magic')
    ,@('Russell Westbrook', 'PG,SG', 'Denver Nuggets')
    ,@('Brandin Podziemski', 'PG,SG', 'Golden State Warriors')
    ,@('Payton Pritchard', 'PG,SG', 'Boston Celtics')
    ,@('Jaylen Brown', 'SG,SF', 'Boston Celtics')
    ,@('Deni Avdija', 'SF,PF', 'Portland Trail Blazers')
    ,@('Pascal Siakam', 'SF,PF,C', 'Indiana Pacers')
    ,@('Zaccharie Risacher', 'SF', 'Atlanta Hawks')
    ,@('Chet Holmgren', 'PF,C', 'Oklahoma City Thunder')
    ,@('Nikola Jokic', 'C', 'Denver Nuggets')
    ,@('Jalen Green', 'PG,SG', 'Houston Rockets')
    ,@('Jerami Grant', 'SF,PF', 'Portland Trail Blazers')
    ,@('Paolo Banchero', 'SF,PF', 'Orlando Magic')
    ,@('Chris Paul', 'PG', 'San Antonio Spurs')
    ,@('Stephon Castle', 'PG,SG', 'San Antonio Spurs')
    ,@('Rudy Gobert', 'C', 'Minnesota Timberwolves')
    ,@('Jalen Suggs', 'PG,SG', 'Orlando Magic')
    ,@('Jakob Poeltl', 'C', 'Toronto Raptors')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
}

Write-Output "Updated roster table rows 2-19"
